$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.199.06"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "1.611.82"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "302.41"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").Value = "0.3784"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "51.79"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "0.3528"
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("D10").Value = "0.08084"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "1.198"
$ws.Range("E11").Value = "  -3.04%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "21.98"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "6.358"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "7.225"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").Value = "0.00001205"
$ws.Range("E16").Value = "  -3.67%  "
$ws.Range("D17").Value = "1.608.35"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "94.20"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "6.507"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "17.14"
$ws.Range("E22").Value = "  -4.47%  "
$ws.Range("D23").Value = "12.31"
$ws.Range("E23").Value = "  -3.57%  "
$ws.Range("D24").Value = "23.180.71"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").Value = "2.504"
$ws.Range("E25").Value = "  +3.31%  "
$ws.Range("D26").Value = "3.013"
$ws.Range("E26").Value = "  -6.91%  "
$ws.Range("D27").Value = "20.83"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").Value = "150.94"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").Value = "5.233"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("D30").Value = "132.27"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").Value = "1.792.58"
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").Value = "1.066"
$ws.Range("E32").Value = "  +10.50%  "
$ws.Range("D33").Value = "6.461"
$ws.Range("E33").Value = "  -5.84%  "
$ws.Range("D34").Value = "2.101"
$ws.Range("E34").Value = "  -8.86%  "
$ws.Range("D35").Value = "11.37"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").Value = "0.02698"
$ws.Range("E36").Value = "  -3.88%  "
$ws.Range("D37").Value = "0.08722"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "0.2442"
$ws.Range("E38").Value = "  -3.95%  "
$ws.Range("D39").Value = "0.06925"
$ws.Range("E39").Value = "  -4.12%  "
$ws.Range("D40").Value = "5.810"
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("D41").Value = "1.321"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("D42").Value = "0.6845"
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").Value = "15.28"
$ws.Range("E44").Value = "  -6.89%  "
$ws.Range("D45").Value = "1.002"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "0.6279"
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").Value = "3.941"
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").Value = "126.94"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").Value = "1.163"
$ws.Range("E51").Value = "  -3.88%  "
